{"js": "// Replace \"del lector de retina\" with \"de reconocimiento facial\" in the\n// \"La base de datos debe comprobar...\" requirement bullet, leaving every\n// other occurrence of \"lector de retina\" elsewhere in the document intact.\n\nconst oldSentence =\n  \"La base de datos debe comprobar, con ayuda del lector de retina, los usuarios registrados para permitir o denegar el ingreso a la instituci\u00f3n.\";\nconst newSentence =\n  \"La base de datos debe comprobar, con ayuda de reconocimiento facial, los usuarios registrados para permitir o denegar el ingreso a la instituci\u00f3n.\";\n\nconst body = context.document.body;\nconst results = body.search(oldSentence, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target sentence not found: \"' + oldSentence + '\"');\n}\n\n// Replace the whole sentence in place (keeps the same run/paragraph, only\n// the run's text content changes, matching the author's edit).\nresults.items[0].insertText(newSentence, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace \"del lector de retina\" with \"de reconocimiento facial\" in the\n# \"La base de datos debe comprobar...\" requirement bullet, leaving every\n# other occurrence of \"lector de retina\" elsewhere in the document intact.\n\n$d = $word.ActiveDocument\n\n$oldSentence = \"La base de datos debe comprobar, con ayuda del lector de retina, los usuarios registrados para permitir o denegar el ingreso a la instituci\u00f3n.\"\n$newSentence = \"La base de datos debe comprobar, con ayuda de reconocimiento facial, los usuarios registrados para permitir o denegar el ingreso a la instituci\u00f3n.\"\n\n$found = $false\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*La base de datos debe comprobar*lector de retina*\") {\n        $rng = $p.Range\n        $rng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)\n        $found = $true\n    }\n}\n\nif (-not $found) {\n    throw \"Target paragraph not found: 'La base de datos debe comprobar...'\"\n}\n"}
